$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'почва'
$ws.Range("C2").Value = 0.3885720750025864
$ws.Range("D2").Value = 'устойчивость почва, почва грунт, режим почва, самоочищение почва'

$ws.Range("B3").Value = 'загрязнение'
$ws.Range("C3").Value = 0.180935045705758
$ws.Range("D3").Value = 'нефтяной загрязнение, источник загрязнение, ореол загрязнение, ореол нефтяной загрязнение'

$ws.Range("B4").Value = 'ландшафт'
$ws.Range("C4").Value = 0.1449324059469544
$ws.Range("D4").Value = 'компонент ландшафт, ландшафтно геохимический, ландшафтно, состояние ландшафт'

$ws.Range("B5").Value = 'геохимический'
$ws.Range("C5").Value = 0.1381817999887385
$ws.Range("D5").Value = 'эколого геохимический, эколого геохимический карта, геохимический поле, ландшафтно геохимический'

$ws.Range("B6").Value = 'поллютант'
$ws.Range("C6").Value = 0.1365767393522176
$ws.Range("D6").Value = 'поллютант природный'

$ws.Range("B7").Value = 'почвенный'
$ws.Range("C7").Value = 0.1084635580499261
$ws.Range("D7").Value = 'почвенный покров, почвенный экосистема'

$ws.Range("B8").Value = 'моск ун тот'
$ws.Range("C8").Value = 0.1037261012262532
$ws.Range("D8").Value = ''

$ws.Range("B9").Value = 'моск ун'
$ws.Range("C9").Value = 0.1037261012262532
$ws.Range("D9").Value = 'моск ун тот, вестн моск ун, во моск ун'

$ws.Range("B10").Value = 'моск'
$ws.Range("C10").Value = 0.09752222461185295
$ws.Range("D10").Value = 'моск ун тот, моск ун, вестн моск, вестн моск ун'

$ws.Range("B11").Value = 'география'
$ws.Range("C11").Value = 0.09589309825489614
$ws.Range("D11").Value = 'сера география, тот сера география'

$ws.Range("B12").Value = 'техногенный'
$ws.Range("C12").Value = 0.0955063826114305
$ws.Range("D12").Value = 'техногенный углеводород, техногенный поток, техногенный поток углеводород, природный техногенный'

$ws.Range("B13").Value = 'тот сера'
$ws.Range("C13").Value = 0.08989595439608609
$ws.Range("D13").Value = 'ун тот сера, тот сера география'

$ws.Range("B14").Value = 'ун тот сера'
$ws.Range("C14").Value = 0.08989595439608609
$ws.Range("D14").Value = ''

$ws.Range("B15").Value = 'пау'
$ws.Range("C15").Value = 0.08989595439608609
$ws.Range("D15").Value = 'ассоциация пау'

$ws.Range("B16").Value = 'сера география'
$ws.Range("C16").Value = 0.08989595439608609
$ws.Range("D16").Value = 'тот сера география'

$ws.Range("B17").Value = 'тот сера география'
$ws.Range("C17").Value = 0.08989595439608609
$ws.Range("D17").Value = ''

$ws.Range("B18").Value = 'самоочищение'
$ws.Range("C18").Value = 0.08970626442700558
$ws.Range("D18").Value = 'потенциал самоочищение, самоочищение почва'

$ws.Range("B19").Value = 'углеводород'
$ws.Range("C19").Value = 0.0866672904902405
$ws.Range("D19").Value = 'техногенный углеводород, углеводородный, полициклический ароматический углеводород, углеводородный геохимический'

$ws.Range("B20").Value = 'карта'
$ws.Range("C20").Value = 0.08666181382866055
$ws.Range("D20").Value = 'эколого геохимический карта, геохимический карта'

$ws.Range("B21").Value = 'природный'
$ws.Range("C21").Value = 0.08623802941483347
$ws.Range("D21").Value = 'природный среда, природный техногенный, нефтепродукт природный, природный техногенный поток'

$ws.Range("B22").Value = 'ун'
$ws.Range("C22").Value = 0.08517805704163238
$ws.Range("D22").Value = 'моск ун тот, моск ун, ун тот сера, ун тот'

$ws.Range("B23").Value = 'ун тот'
$ws.Range("C23").Value = 0.08517805704163238
$ws.Range("D23").Value = 'моск ун тот, ун тот сера'

$ws.Range("B24").Value = 'тот'
$ws.Range("C24").Value = 0.08313586864282774
$ws.Range("D24").Value = 'моск ун тот, тот сера, ун тот сера, тот сера география'

$ws.Range("B25").Value = 'эколого геохимический'
$ws.Range("C25").Value = 0.08298088098100256
$ws.Range("D25").Value = 'эколого геохимический карта, прогнозный эколого геохимический, эколого геохимический состояние'

$ws.Range("B26").Value = 'нефтепродукт'
$ws.Range("C26").Value = 0.07969188806726681
$ws.Range("D26").Value = 'нефть нефтепродукт, аккумуляция нефть нефтепродукт, нефтепродукт природный, содержание нефтепродукт'

$ws.Range("B27").Value = 'нефть'
$ws.Range("C27").Value = 0.07876325837688636
$ws.Range("D27").Value = 'нефть нефтепродукт, аккумуляция нефть, аккумуляция нефть нефтепродукт, нефть окружающий среда'

$ws.Range("B28").Value = 'пиковский'
$ws.Range("C28").Value = 0.07801777968948236
$ws.Range("D28").Value = 'геннадий пиковский'

$ws.Range("B29").Value = 'ореол'
$ws.Range("C29").Value = 0.0777454291700715
$ws.Range("D29").Value = 'ореол загрязнение, ореол нефтяной загрязнение, ореол нефтяной'

$ws.Range("B30").Value = 'вестн моск'
$ws.Range("C30").Value = 0.07606580756591901
$ws.Range("D30").Value = 'вестн моск ун'

$ws.Range("B31").Value = 'вестн моск ун'
$ws.Range("C31").Value = 0.07606580756591901
$ws.Range("D31").Value = ''
